$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("AH2").Value = 12
$ws.Range("AD3").Value = 501
$ws.Range("AE3").Value = 6.5
$ws.Range("I3").Value = 2.6
$ws.Range("R3").Value = 2.05
$ws.Range("S3").Value = 1.7
$ws.Range("V3").Value = 12
$ws.Range("Z3").Value = 7
$ws.Range("J4").Value = 1.1
$ws.Range("K4").Value = 7
$ws.Range("J5").Value = 1.18
$ws.Range("K5").Value = 4.5
$ws.Range("AI6").Value = 34
$ws.Range("I6").Value = 3.1
$ws.Range("J6").Value = 1.11
$ws.Range("K6").Value = 6.5
$ws.Range("U6").Value = 10
$ws.Range("W6").Value = 23
$ws.Range("AF7").Value = 19
$ws.Range("G7").Value = 2.1
$ws.Range("I7").Value = 4.1
$ws.Range("N7").Value = 2.7
$ws.Range("O7").Value = 1.44
$ws.Range("G12").Value = 1.75
$ws.Range("I12").Value = 4.5
$ws.Range("N12").Value = 2
$ws.Range("U12").Value = 8.5
$ws.Range("W12").Value = 15
$ws.Range("K13").Value = 9
$ws.Range("N13").Value = 2.1
$ws.Range("O13").Value = 1.7
$ws.Range("AE14").Value = 9.5
$ws.Range("AF14").Value = 17
$ws.Range("AG14").Value = 13
$ws.Range("G14").Value = 2.2
$ws.Range("H14").Value = 3.1
$ws.Range("I14").Value = 3.5
$ws.Range("U14").Value = 10
$ws.Range("X14").Value = 19
$ws.Range("L15").Value = 1.3
$ws.Range("M15").Value = 3.4
$ws.Range("N15").Value = 2
$ws.Range("U15").Value = 29
$ws.Range("AE17").Value = 6.4
$ws.Range("AF17").Value = 8
$ws.Range("R17").Value = 1.82
$ws.Range("T17").Value = 11
$ws.Range("Y17").Value = 50
$ws.Range("Z17").Value = 8.5
$ws.Range("AE19").Value = 9.5
$ws.Range("AF19").Value = 15
$ws.Range("G19").Value = 2.5
$ws.Range("Z19").Value = 10
$ws.Range("AF22").Value = 13
$ws.Range("AG22").Value = 9.5
$ws.Range("AI22").Value = 17
$ws.Range("G22").Value = 2.55
$ws.Range("I22").Value = 2.3
$ws.Range("J22").Value = 1.03
$ws.Range("L22").Value = 1.2
$ws.Range("T22").Value = 12
$ws.Range("W22").Value = 29
$ws.Range("J23").Value = 1.07
$ws.Range("L23").Value = 1.36
$ws.Range("N27").Value = 2
$ws.Range("N29").Value = 2.08
$ws.Range("O29").Value = 1.73
$ws.Range("AA30").Value = 8.5
$ws.Range("AB30").Value = 19
$ws.Range("AD30").Value = 301
$ws.Range("AE30").Value = 19
$ws.Range("AF30").Value = 41
$ws.Range("AG30").Value = 21
$ws.Range("AH30").Value = 81
$ws.Range("AI30").Value = 51
$ws.Range("AJ30").Value = 51
$ws.Range("G30").Value = 1.42
$ws.Range("H30").Value = 4.33
$ws.Range("I30").Value = 7.5
$ws.Range("N30").Value = 1.7
$ws.Range("O30").Value = 2.1
$ws.Range("R30").Value = 1.95
$ws.Range("S30").Value = 1.8
$ws.Range("T30").Value = 7
$ws.Range("U30").Value = 7
$ws.Range("W30").Value = 9.5
$ws.Range("Y30").Value = 26
$ws.Range("AA31").Value = 7
$ws.Range("AB31").Value = 17
$ws.Range("AD31").Value = 301
$ws.Range("AE31").Value = 15
$ws.Range("AG31").Value = 19
$ws.Range("AH31").Value = 67
$ws.Range("AJ31").Value = 51
$ws.Range("H31").Value = 3.5
$ws.Range("I31").Value = 5.75
$ws.Range("L31").Value = 1.25
$ws.Range("M31").Value = 3.75
$ws.Range("N31").Value = 1.93
$ws.Range("O31").Value = 1.93
$ws.Range("U31").Value = 7.5
$ws.Range("G32").Value = 1.91
$ws.Range("I32").Value = 4
$ws.Range("U32").Value = 9
$ws.Range("W32").Value = 17
$ws.Range("AE33").Value = 10.75
$ws.Range("J33").Value = 1.06
$ws.Range("M33").Value = 3.25
$ws.Range("N33").Value = 1.87
$ws.Range("O33").Value = 1.83
$ws.Range("P33").Value = 1.42
$ws.Range("Q33").Value = 2.65
$ws.Range("R33").Value = 1.7
$ws.Range("U33").Value = 9.75
$ws.Range("AA34").Value = 7.6
$ws.Range("AE34").Value = 5.3
$ws.Range("AF34").Value = 6.1
$ws.Range("AG34").Value = 8.75
$ws.Range("AH34").Value = 10
$ws.Range("AI34").Value = 14.5
$ws.Range("AJ34").Value = 40
$ws.Range("G34").Value = 5.8
$ws.Range("H34").Value = 3.7
$ws.Range("I34").Value = 1.52
$ws.Range("J34").Value = 1.08
$ws.Range("K34").Value = 6.6
$ws.Range("L34").Value = 1.38
$ws.Range("M34").Value = 2.8
$ws.Range("N34").Value = 2.12
$ws.Range("O34").Value = 1.65
$ws.Range("P34").Value = 1.42
$ws.Range("Q34").Value = 2.67
$ws.Range("R34").Value = 2.25
$ws.Range("S34").Value = 1.57
$ws.Range("T34").Value = 12.5
$ws.Range("U34").Value = 35
$ws.Range("V34").Value = 20
$ws.Range("W34").Value = 120
$ws.Range("Z34").Value = 6.6
$ws.Range("AA35").Value = 6.5
$ws.Range("AE35").Value = 12
$ws.Range("H35").Value = 3.25
$ws.Range("L35").Value = 1.26
$ws.Range("X35").Value = 14.5
$ws.Range("L36").Value = 1.23
$ws.Range("AG37").Value = 26
$ws.Range("G37").Value = 1.3
$ws.Range("H37").Value = 5.25
$ws.Range("I37").Value = 9
$ws.Range("J37").Value = 1.05
$ws.Range("K37").Value = 11
$ws.Range("N38").Value = 1.83
$ws.Range("O38").Value = 2.03
$ws.Range("N39").Value = 1.83
$ws.Range("O39").Value = 2.03
$ws.Range("N40").Value = 1.67
$ws.Range("O40").Value = 2.15
$ws.Range("AA41").Value = 6.1
$ws.Range("AB41").Value = 25
$ws.Range("AC41").Value = 200
$ws.Range("AE41").Value = 6.7
$ws.Range("AF41").Value = 16.5
$ws.Range("AG41").Value = 15
$ws.Range("AJ41").Value = 90
$ws.Range("G41").Value = 2.18
$ws.Range("H41").Value = 2.82
$ws.Range("J41").Value = 1.15
$ws.Range("K41").Value = 4.2
$ws.Range("L41").Value = 1.65
$ws.Range("M41").Value = 1.98
$ws.Range("N41").Value = 2.87
$ws.Range("O41").Value = 1.31
$ws.Range("P41").Value = 1.65
$ws.Range("Q41").Value = 1.98
$ws.Range("R41").Value = 2.37
$ws.Range("S41").Value = 1.45
$ws.Range("T41").Value = 4.85
$ws.Range("V41").Value = 10.5
$ws.Range("X41").Value = 27
$ws.Range("Y41").Value = 60
$ws.Range("Z41").Value = 4.45
$ws.Range("AA42").Value = 5.9
$ws.Range("AB42").Value = 19
$ws.Range("AC42").Value = 150
$ws.Range("AE42").Value = 7.8
$ws.Range("AH42").Value = 50
$ws.Range("G42").Value = 2.15
$ws.Range("H42").Value = 2.95
$ws.Range("I42").Value = 3.5
$ws.Range("J42").Value = 1.11
$ws.Range("K42").Value = 6
$ws.Range("L42").Value = 1.5
$ws.Range("M42").Value = 2.25
$ws.Range("N42").Value = 2.45
$ws.Range("O42").Value = 1.42
$ws.Range("P42").Value = 1.53
$ws.Range("Q42").Value = 2.18
$ws.Range("R42").Value = 2.07
$ws.Range("S42").Value = 1.6
$ws.Range("T42").Value = 5.4
$ws.Range("V42").Value = 9.5
$ws.Range("X42").Value = 22
$ws.Range("Y42").Value = 45
$ws.Range("Z42").Value = 6.3
$ws.Range("AB43").Value = 13
$ws.Range("AC43").Value = 41
$ws.Range("AD43").Value = 201
$ws.Range("AE43").Value = 8
$ws.Range("AF43").Value = 10
$ws.Range("AI43").Value = 15
$ws.Range("AJ43").Value = 26
$ws.Range("G43").Value = 3.5
$ws.Range("I43").Value = 2
$ws.Range("N43").Value = 1.93
$ws.Range("O43").Value = 1.93
$ws.Range("P43").Value = 1.36
$ws.Range("Q43").Value = 3
$ws.Range("R43").Value = 1.73
$ws.Range("S43").Value = 2
$ws.Range("T43").Value = 12
$ws.Range("Y43").Value = 34
$ws.Range("Z43").Value = 11
